$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7042609453201294
$ws.Range("B1").Value = 1.236676692962646
$ws.Range("C1").Value = 3.604037523269653
$ws.Range("D1").Value = 2.804191827774048
$ws.Range("E1").Value = 1.51946747303009
